$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NOFOpt")
$ws.Activate()

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 100
$ws.Range("B4").Value = 300
$ws.Range("C4").Value = 300

$ws.Range("A5:C5").Delete()

$ws.Range("G10").Select()
